$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (45406 -> 45436, i.e. one month later)
$ws.Range("A1").Value = 45436

# Update the prices in column D
$ws.Range("D14").Value = 1266.597
$ws.Range("D15").Value = 1546.566
$ws.Range("D16").Value = 1817.002
